$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that currently sits right after
#    "MP73010" in the title line.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the paragraph that contains "Ben changing things up!" and
#    locate the (empty) paragraph immediately following it - that is
#    where the new paragraph needs to be inserted.
# ------------------------------------------------------------------
$benIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Ben changing things up!*") {
        $benIndex = $i
        break
    }
}

$targetIndex = $benIndex + 1
$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $targetPara.Range

# Insert a brand-new paragraph straight after it; this is where the
# "Jianqi changing things up!" text (and the relocated bookmark) will
# live.
$targetRange.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range

# ------------------------------------------------------------------
# 3. Populate the new paragraph with the exact run/proofErr structure
#    from the edit: "J" (hinted East-Asian run) + "ianqi" + the rest
#    of the sentence, wrapped in spell-check proofErr markers, plus
#    the eastAsia language tagging on the paragraph mark itself, and
#    the "_GoBack" bookmark at the very end of the paragraph.
# ------------------------------------------------------------------
$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>J</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>ianqi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> changing things up! </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newRange.InsertXML($newParaXml)
